$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds text-formatted numbers (e.g. "63.263.06"); keep them as text
# rather than letting Excel auto-convert numeric-looking strings to numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "63.263.06"
$ws.Range("E2").Value = "  -1.25%  "
$ws.Range("D3").Value = "2.693.78"
$ws.Range("E3").Value = "  -2.49%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "556.58"
$ws.Range("E5").Value = "  -3.54%  "
$ws.Range("D6").Value = "156.97"
$ws.Range("E6").Value = "  -1.58%  "
$ws.Range("E7").Value = "  +0.10%  "
$ws.Range("D8").Value = "0.581"
$ws.Range("E8").Value = "  -3.32%  "
$ws.Range("D9").Value = "0.105"
$ws.Range("E9").Value = "  -4.23%  "
$ws.Range("E10").Value = "  +0.13%  "
$ws.Range("D11").Value = "0.368"
$ws.Range("E11").Value = "  -4.51%  "
$ws.Range("D12").Value = "5.42"
$ws.Range("E12").Value = "  -6.60%  "
$ws.Range("D13").Value = "3.174.64"
$ws.Range("E13").Value = "  -2.34%  "
$ws.Range("D14").Value = "26.33"
$ws.Range("E14").Value = "  -2.35%  "
$ws.Range("D15").Value = "63.126.28"
$ws.Range("E15").Value = "  -0.84%  "
$ws.Range("D16").Value = "0.0000144"
$ws.Range("D17").Value = "2.705.46"
$ws.Range("E17").Value = "  -2.26%  "
$ws.Range("D18").Value = "12.04"
$ws.Range("E18").Value = "  -0.62%  "
$ws.Range("E19").Value = "  -5.45%  "
$ws.Range("D20").Value = "343.48"
$ws.Range("E20").Value = "  -4.63%  "
$ws.Range("E21").Value = "  -4.95%  "
$ws.Range("D22").Value = "0.997"
$ws.Range("E22").Value = "  -0.28%  "
$ws.Range("E23").Value = "  -4.30%  "
$ws.Range("D24").Value = "63.72"
$ws.Range("E24").Value = "  -2.00%  "
$ws.Range("E25").Value = "  -0.65%  "
$ws.Range("D26").Value = "0.998"
$ws.Range("E26").Value = "  -0.13%  "
$ws.Range("D27").Value = "8.05"
$ws.Range("E27").Value = "  -5.66%  "
$ws.Range("D28").Value = "0.0₃0860"
$ws.Range("E28").Value = "  -5.07%  "
$ws.Range("E29").Value = "  -1.45%  "
$ws.Range("E30").Value = "  +5.22%  "
$ws.Range("D31").Value = "7.05"
$ws.Range("E31").Value = "  -3.49%  "
$ws.Range("D32").Value = "165.16"
$ws.Range("E32").Value = "  -3.09%  "
$ws.Range("E33").Value = "  -0.01%  "
$ws.Range("D34").Value = "19.53"
$ws.Range("E34").Value = "  -3.20%  "
$ws.Range("D35").Value = "4.76"
$ws.Range("E35").Value = "  -3.38%  "
$ws.Range("E36").Value = "  -2.81%  "
$ws.Range("D37").Value = "1.78"
$ws.Range("E37").Value = "  -1.56%  "
$ws.Range("D38").Value = "340.00"
$ws.Range("E38").Value = "  -2.54%  "
$ws.Range("E39").Value = "  -5.56%  "
$ws.Range("E40").Value = "  -3.59%  "
$ws.Range("E41").Value = "  -5.33%  "
$ws.Range("D42").Value = "38.08"
$ws.Range("E42").Value = "  -2.57%  "
$ws.Range("D43").Value = "20.88"
$ws.Range("E43").Value = "  -4.12%  "
$ws.Range("D44").Value = "20.27"
$ws.Range("E44").Value = "  -5.52%  "
$ws.Range("E45").Value = "  -1.87%  "
$ws.Range("E46").Value = "  -4.24%  "
$ws.Range("D47").Value = "0.998"
$ws.Range("E47").Value = "  -0.03%  "
$ws.Range("B48").Value = "WhiteBITCoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D48").Value = "11.07"
$ws.Range("E48").Value = "  +0.10%  "
$ws.Range("B49").Value = "Aave"
$ws.Range("C49").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D49").Value = "130.47"
$ws.Range("E49").Value = "  -5.55%  "
$ws.Range("D50").Value = "0.0974"
$ws.Range("E50").Value = "  -3.68%  "
$ws.Range("D51").Value = "2.104.42"
$ws.Range("E51").Value = "  -0.82%  "
